# Correct the spelling of several governorate names (same governorates,
# just corrected/standardized Arabic spelling). This changes the text of
# the corresponding shared-string entries; row order in the sheet is
# untouched. Applied in the order that reproduces the target shared
# string table ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "الأقصر"
$ws.Range("A6").Value  = "الإسكندرية"
$ws.Range("A2").Value  = "القليوبية"
$ws.Range("A4").Value  = "الشرقية"
$ws.Range("A28").Value = "أسوان"
$ws.Range("A17").Value = "بني سويف"
$ws.Range("A26").Value = "الوادي الجديد"

# Update the remembered cell selection on the sheet.
$null = $ws.Range("D26").Select()
